$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force the cell to stay a text value (avoid Excel auto-converting
    # numeric-looking strings like "1.004" into numbers), then restore
    # the default "Normal" style so no stray formatting is introduced.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell $ws 'D2' '27.882.66'
Set-TextCell $ws 'E2' '  -0.64%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell $ws 'D3' '1.908.10'
Set-TextCell $ws 'E3' '  -0.17%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell $ws 'D4' '1.004'
Set-TextCell $ws 'E4' '  -0.16%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell $ws 'D5' '314.11'
Set-TextCell $ws 'E5' '  -0.91%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell $ws 'D6' '1.003'
Set-TextCell $ws 'E6' '  -0.17%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell $ws 'D7' '0.5006'
Set-TextCell $ws 'E7' '  +3.82%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell $ws 'D8' '0.3818'
Set-TextCell $ws 'E8' '  +0.06%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws 'D9' '0.07293'
Set-TextCell $ws 'E9' '  -0.98%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell $ws 'D10' '0.9098'
Set-TextCell $ws 'E10' '  -2.66%  '

$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell $ws 'D11' '20.89'
Set-TextCell $ws 'E11' '  +0.23%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 'D12' '1.936.51'
Set-TextCell $ws 'E12' '  +1.45%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D13' '0.07681'
Set-TextCell $ws 'E13' '  -1.36%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws 'D14' '5.480'
Set-TextCell $ws 'E14' '  -0.61%  '

$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws 'D15' '91.85'
Set-TextCell $ws 'E15' '  +0.03%  '

$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell $ws 'D16' '1.005'
Set-TextCell $ws 'E16' '  -0.14%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D17' '0.000008730'
Set-TextCell $ws 'E17' '  -1.25%  '

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D18' '1.003'
Set-TextCell $ws 'E18' '  -0.15%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 'D19' '27.931.80'
Set-TextCell $ws 'E19' '  -0.55%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws 'D20' '14.60'
Set-TextCell $ws 'E20' '  -1.67%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws 'D21' '5.174'
Set-TextCell $ws 'E21' '  -0.17%  '

$ws.Range('B22').Value = 'Cosmos'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D22' '10.83'
Set-TextCell $ws 'E22' '  -0.95%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D23' '6.586'
Set-TextCell $ws 'E23' '  -0.69%  '

$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D24' '154.35'
Set-TextCell $ws 'E24' '  -0.85%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D25' '1.881'
Set-TextCell $ws 'E25' '  -2.09%  '

$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 'D26' '2.221'
Set-TextCell $ws 'E26' '  +4.62%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D27' '18.42'
Set-TextCell $ws 'E27' '  -0.76%  '

$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws 'D28' '115.41'
Set-TextCell $ws 'E28' '  -1.13%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D29' '4.920'
Set-TextCell $ws 'E29' '  -0.86%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D30' '0.08991'
Set-TextCell $ws 'E30' '  +0.30%  '

$ws.Range('B31').Value = 'HuobiToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws 'D31' '3.209'
Set-TextCell $ws 'E31' '  -3.10%  '

$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws 'D32' '1.232'
Set-TextCell $ws 'E32' '  -1.87%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D33' '0.7653'
Set-TextCell $ws 'E33' '  -1.71%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D34' '4.656'
Set-TextCell $ws 'E34' '  -0.60%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D35' '0.02064'
Set-TextCell $ws 'E35' '  +0.37%  '

$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws 'D36' '2.543'
Set-TextCell $ws 'E36' '  -4.47%  '

$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws 'D37' '0.5580'
Set-TextCell $ws 'E37' '  +1.65%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws 'D38' '1.093'
Set-TextCell $ws 'E38' '  -1.66%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws 'D39' '3.024'
Set-TextCell $ws 'E39' '  +1.17%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D40' '0.05253'
Set-TextCell $ws 'E40' '  -1.31%  '

$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws 'D41' '6.946'
Set-TextCell $ws 'E41' '  -1.12%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D42' '8.508'
Set-TextCell $ws 'E42' '  -0.36%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws 'D43' '0.1511'
Set-TextCell $ws 'E43' '  -1.20%  '

$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws 'D44' '111.13'
Set-TextCell $ws 'E44' '  +2.37%  '

$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell $ws 'D45' '0.4833'
Set-TextCell $ws 'E45' '  -0.21%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D46' '10.58'
Set-TextCell $ws 'E46' '  -1.48%  '

$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws 'D47' '1.003'
Set-TextCell $ws 'E47' '  -0.24%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D48' '1.630'
Set-TextCell $ws 'E48' '  -1.66%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D49' '67.57'
Set-TextCell $ws 'E49' '  -0.81%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws 'D50' '0.06066'
Set-TextCell $ws 'E50' '  -0.37%  '

$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell $ws 'D51' '0.9024'
Set-TextCell $ws 'E51' '  +0.19%  '
